$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data held in row 2 and row 3 for the columns that differ
# per-record (date, volume, prices, unit description, $/kg, kg/unit).
$cols = @("D", "M", "N", "O", "P", "Q", "S", "T")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value2 = $cell3.Value2
    $cell3.Value2 = $tmp
}
